$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.116910696029663
$ws.Range("B1").Value = 2.258161067962646
$ws.Range("C1").Value = 10.4215784072876
$ws.Range("D1").Value = 1.711314082145691
$ws.Range("E1").Value = 1.29023003578186
